$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Add the new row of data (row 30) for the "Brutto årslønn" variable
$ws.Range("A30").Value = "Brutto årslønn"
$ws.Range("B30").Value = "brutto_arslonn_vasket"
$ws.Range("C30").Value = "snitt_as_num_single"
$ws.Range("D30").Value = "Brutto årslønn"
$ws.Range("E30").Value = "Vi har tatt bort svar som er under 300 000 og over 1 000 000 kr, og gjennomsnittet inkluderer bare de som oppgir å arbeide med det de er utdannet til."

# Resize the Excel table (ListObject) so the new row becomes part of Table1
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:E30"))

# The data validation on column D previously spanned the whole column
# (D1:D1048576). Row 30 must be excluded from it, splitting it into
# D1:D29 and D31:D1048576 while keeping the rest of the column validated.
$ws.Range("D30").Validation.Delete()

# Update the selection/view to point at the newly added row
$ws.Range("A30:E30").Select()
